$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.286.40'
$ws.Range("E2").Value = '  +6.54%  '

# Row 3
$ws.Range("D3").Value = '3.546.71'
$ws.Range("E3").Value = '  +10.07%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.37'
$ws.Range("E5").Value = '  +9.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '552.41'
$ws.Range("E6").Value = '  +4.28%  '

# Row 7
$ws.Range("D7").Value = '3.536.35'
$ws.Range("E7").Value = '  +9.93%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.608'
$ws.Range("E8").Value = '  +2.36%  '

# Row 9
$ws.Range("E9").Value = '  -0.19%  '

# Row 10
$ws.Range("E10").Value = '  +4.33%  '

# Row 11
$ws.Range("E11").Value = '  +15.59%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.89'
$ws.Range("E12").Value = '  +2.52%  '

# Row 13
$ws.Range("E13").Value = '  +6.77%  '

# Row 14
$ws.Range("E14").Value = '  +3.16%  '

# Row 15
$ws.Range("D15").Value = '4.104.34'
$ws.Range("E15").Value = '  +9.69%  '

# Row 16
$ws.Range("D16").Value = '3.544.23'
$ws.Range("E16").Value = '  +10.02%  '

# Row 17
$ws.Range("E17").Value = '  +4.79%  '

# Row 18
$ws.Range("D18").Value = '67.293.79'
$ws.Range("E18").Value = '  +6.83%  '

# Row 19
$ws.Range("E19").Value = '  +5.59%  '

# Row 20
$ws.Range("E20").Value = '  +7.92%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.994'
$ws.Range("E21").Value = '  +2.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '428.77'
$ws.Range("E22").Value = '  +16.94%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '85.50'
$ws.Range("E23").Value = '  +5.27%  '

# Row 24
$ws.Range("E24").Value = '  +3.60%  '

# Row 25
$ws.Range("E25").Value = '  +5.19%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.11'
$ws.Range("E26").Value = '  +0.59%  '

# Row 27
$ws.Range("E27").Value = '  +9.54%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.09'
$ws.Range("E28").Value = '  +6.46%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.00'
$ws.Range("E29").Value = '  +9.97%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.38'
$ws.Range("E30").Value = '  +6.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '645.05'
$ws.Range("E31").Value = '  +1.39%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.68'
$ws.Range("E32").Value = '  +3.02%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.75'
$ws.Range("E33").Value = '  +4.12%  '

# Row 34
$ws.Range("E34").Value = '  +4.53%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.76'
$ws.Range("E35").Value = '  +5.23%  '

# Row 36
$ws.Range("D36").Value = '0.0₃0830'
$ws.Range("E36").Value = '  +15.95%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.61'
$ws.Range("E37").Value = '  +4.64%  '

# Row 38
$ws.Range("E38").Value = '  +18.79%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.05%  '

# Row 40
$ws.Range("E40").Value = '  +4.04%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.34'
$ws.Range("E41").Value = '  +14.25%  '

# Row 42
$ws.Range("E42").Value = '  +0.26%  '

# Row 43
$ws.Range("D43").Value = '3.041.72'
$ws.Range("E43").Value = '  +5.22%  '

# Row 44
$ws.Range("E44").Value = '  +5.37%  '

# Row 45
$ws.Range("E45").Value = '  +10.95%  '

# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.37'
$ws.Range("E46").Value = '  +12.59%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.83'
$ws.Range("E47").Value = '  +5.55%  '

# Row 48
$ws.Range("E48").Value = '  +5.89%  '

# Row 49
$ws.Range("E49").Value = '  +5.06%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.70'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '141.13'
$ws.Range("E51").Value = '  +4.22%  '

